# Appium test-data workbook update
# - happy_path_data!A2: new valid-user email (keeps the "quote-prefix" text style)
# - invalid_user: new "password" column (B) with header + sample value, and the
#   existing sample email in A2 is replaced with a different bogus address
# - selection / active-sheet bookkeeping to match the authored edit

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("happy_path_data")
$ws2 = $wb.Worksheets.Item("invalid_user")

# --- happy_path_data ---------------------------------------------------
# Leading apostrophe re-asserts the text ("quote prefix") formatting that the
# cell already had, instead of resetting it to the plain default style.
$ws1.Range("A2").Value = "'validuser2233@email.com"

# --- invalid_user --------------------------------------------------------
# Replace the sample (invalid) email with a different one.
$ws2.Range("A2").Value = "'animesh5678@gmail.com"

# Add the new "password" column next to "email".
$ws2.Range("B1").Value = "password"
$ws2.Range("A1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)   # xlPasteFormats - match the header style

$ws2.Range("B2").Value = "'Welcome@1"
$ws2.Range("A2").Copy()
$ws2.Range("B2").PasteSpecial(-4122)   # xlPasteFormats - match the quote-prefixed style

$excel.CutCopyMode = 0

# --- selection / active tab bookkeeping -----------------------------------
# Select on the sheet that should end up *not* active first, then finish on
# the sheet that should be the active / visible tab.
$ws2.Activate()
$ws2.Range("C4").Select()

$ws1.Activate()
$ws1.Range("B10").Select()
